$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Actual Consumption (MW)" / "Timestamp" data that replaces the
# previous portfolio's rows 2-33 (old rows 34-51 are removed entirely).
$data = @(
    @(5152, 45740),
    @(5127, 45740.01041666666),
    @(5090, 45740.02083333334),
    @(5049, 45740.03125),
    @(5035, 45740.04166666666),
    @(5015, 45740.05208333334),
    @(4989, 45740.0625),
    @(4956, 45740.07291666666),
    @(4965, 45740.08333333334),
    @(4972, 45740.09375),
    @(4968, 45740.10416666666),
    @(4985, 45740.11458333334),
    @(5001, 45740.125),
    @(5022, 45740.13541666666),
    @(5059, 45740.14583333334),
    @(5103, 45740.15625),
    @(5246, 45740.16666666666),
    @(5333, 45740.17708333334),
    @(5420, 45740.1875),
    @(5528, 45740.19791666666),
    @(5758, 45740.20833333334),
    @(5863, 45740.21875),
    @(6016, 45740.22916666666),
    @(6211, 45740.23958333334),
    @(6501, 45740.25),
    @(6683, 45740.26041666666),
    @(6812, 45740.27083333334),
    @(6886, 45740.28125),
    @(7031, 45740.29166666666),
    @(7044, 45740.30208333334),
    @(7069, 45740.3125),
    @(6997, 45740.32291666666)
)

# First remove the now-obsolete rows (old rows 34-51) by deleting entire
# rows so everything below shifts up and the sheet dimension shrinks.
$ws.Range("A34:B51").EntireRow.Delete()

# Overwrite rows 2-33 with the refreshed data values.
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}
